$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.847.58"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "2.603.31"
$ws.Range("E3").Value = "  -1.52%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.88"
$ws.Range("E5").Value = "  -1.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.63"
$ws.Range("E6").Value = "  -0.71%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.527"
$ws.Range("E8").Value = "  -3.53%  "
$ws.Range("D9").Value = "2.602.26"
$ws.Range("E9").Value = "  -1.53%  "
$ws.Range("E10").Value = "  -3.02%  "
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.365"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("E13").Value = "  -1.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.07"
$ws.Range("E14").Value = "  -3.16%  "
$ws.Range("D15").Value = "3.079.26"
$ws.Range("E15").Value = "  -1.53%  "
$ws.Range("E16").Value = "  -2.59%  "
$ws.Range("D17").Value = "67.191.74"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").Value = "2.592.89"
$ws.Range("E18").Value = "  -2.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.61"
$ws.Range("E19").Value = "  -4.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.78"
$ws.Range("E20").Value = "  -4.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "352.98"
$ws.Range("E21").Value = "  -2.39%  "
$ws.Range("E22").Value = "  -3.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.61"
$ws.Range("E23").Value = "  -3.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.44"
$ws.Range("E24").Value = "  -5.19%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("E26").Value = "  -5.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "69.10"
$ws.Range("E27").Value = "  -2.34%  "
$ws.Range("D28").Value = "2.741.09"
$ws.Range("E28").Value = "  -1.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").Value = "0.0₃0985"
$ws.Range("E30").Value = "  -3.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "538.27"
$ws.Range("E31").Value = "  -3.31%  "
$ws.Range("E32").Value = "  +1.28%  "
$ws.Range("E33").Value = "  -4.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.85"
$ws.Range("E34").Value = "  -3.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.133"
$ws.Range("E35").Value = "  -0.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("E37").Value = "  -4.92%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.92"
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.81"
$ws.Range("E39").Value = "  -2.31%  "
$ws.Range("E40").Value = "  -2.63%  "
$ws.Range("E41").Value = "  +1.91%  "
$ws.Range("E42").Value = "  -1.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.10"
$ws.Range("E43").Value = "  -3.62%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.40"
$ws.Range("E45").Value = "  -5.11%  "
$ws.Range("E46").Value = "  -1.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "150.06"
$ws.Range("E47").Value = "  -2.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.572"
$ws.Range("E48").Value = "  -3.47%  "
$ws.Range("E49").Value = "  -3.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.70"
$ws.Range("E50").Value = "  -1.97%  "
$ws.Range("E51").Value = "  -1.38%  "
